# "End of factorisation": the per-row "group" column duplicated in column C
# of both sheets only belongs on sheet1's companion sheet. Finish the
# factorisation by removing it from "sheet1" and inserting it as column C
# of "Feuille2" (where the lookup/value columns it sits next to already
# live), instead of keeping a redundant copy on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("sheet1")
$ws2 = $wb.Worksheets.Item("Feuille2")

# Capture the "group" labels currently in sheet1!C2:C15 before removing
# that column.
$groups = @()
for ($r = 2; $r -le 15; $r++) {
    $groups += $ws1.Cells.Item($r, 3).Value2
}

# Remove the now-redundant column C ("group…") from sheet1. Everything to
# the right (the D formula column + the G value column) shifts left by
# one column; the formula referencing the value column auto-adjusts.
$ws1.Range("C1").EntireColumn.Delete()

# Insert a fresh column C on Feuille2 to host that same "group" data (the
# existing C/D/value columns there shift right by one; formulas
# auto-adjust).
$ws2.Range("C1").EntireColumn.Insert()

# Repopulate the newly inserted column with the captured group labels.
for ($r = 2; $r -le 15; $r++) {
    $ws2.Cells.Item($r, 3).Value = $groups[$r - 2]
}

# The bottom two rows sit in a highlighted row band; re-apply the plain
# (non-highlighted) cell style used by the rest of the inserted column so
# the new cells don't pick up that highlight by inheritance.
$ws2.Range("C14:C15").Style = $ws2.Range("D2").Style

# Leave the selection on Feuille2 parked below the data, then switch back
# to sheet1 as the active/visible tab.
$ws2.Range("C17").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("C1").Select() | Out-Null
